$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format cells whose new numeric-looking price strings must stay text
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D12", "D13", "D14", "D15", "D17", "D18", "D19", "D20", "D21", "D22", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '28.710.48'
$ws.Range("E2").Value = '  +1.37%  '
$ws.Range("D3").Value = '1.807.10'
$ws.Range("E3").Value = '  -0.22%  '
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.37%  '
$ws.Range("D5").Value = '328.29'
$ws.Range("E5").Value = '  -2.56%  '
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.49%  '
$ws.Range("D7").Value = '0.4371'
$ws.Range("E7").Value = '  -1.44%  '
$ws.Range("D8").Value = '0.3756'
$ws.Range("E8").Value = '  +6.42%  '
$ws.Range("D9").Value = '46.35'
$ws.Range("E9").Value = '  +1.40%  '
$ws.Range("D10").Value = '0.07662'
$ws.Range("E10").Value = '  +2.65%  '
$ws.Range("E11").Value = '  -1.28%  '
$ws.Range("D12").Value = '22.76'
$ws.Range("E12").Value = '  -0.83%  '
$ws.Range("D13").Value = '1.002'
$ws.Range("E13").Value = '  +0.27%  '
$ws.Range("D14").Value = '6.261'
$ws.Range("E14").Value = '  -0.35%  '
$ws.Range("D15").Value = '7.502'
$ws.Range("E15").Value = '  +2.96%  '
$ws.Range("D16").Value = '1.804.66'
$ws.Range("E16").Value = '  -0.71%  '
$ws.Range("D17").Value = '0.00001091'
$ws.Range("D18").Value = '0.06710'
$ws.Range("E18").Value = '  +0.79%  '
$ws.Range("D19").Value = '81.18'
$ws.Range("E19").Value = '  -0.99%  '
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  +0.44%  '
$ws.Range("D21").Value = '17.63'
$ws.Range("E21").Value = '  +1.76%  '
$ws.Range("D22").Value = '6.273'
$ws.Range("E22").Value = '  -2.61%  '
$ws.Range("D23").Value = '28.712.89'
$ws.Range("E23").Value = '  +1.29%  '
$ws.Range("E24").Value = '  -2.15%  '
$ws.Range("D25").Value = '2.456'
$ws.Range("E25").Value = '  +2.74%  '
$ws.Range("D26").Value = '20.54'
$ws.Range("E26").Value = '  -0.96%  '
$ws.Range("D27").Value = '154.95'
$ws.Range("E27").Value = '  -0.14%  '
$ws.Range("D28").Value = '2.359'
$ws.Range("E28").Value = '  -3.93%  '
$ws.Range("D29").Value = '2.014.94'
$ws.Range("E29").Value = '  -0.37%  '
$ws.Range("D30").Value = '1.307'
$ws.Range("E30").Value = '  -0.47%  '
$ws.Range("D31").Value = '130.96'
$ws.Range("E31").Value = '  -1.36%  '
$ws.Range("D32").Value = '3.972'
$ws.Range("E32").Value = '  -2.23%  '
$ws.Range("D33").Value = '5.819'
$ws.Range("E33").Value = '  -2.68%  '
$ws.Range("D34").Value = '0.09193'
$ws.Range("E34").Value = '  -1.87%  '
$ws.Range("D35").Value = '0.2225'
$ws.Range("E35").Value = '  +2.75%  '
$ws.Range("D36").Value = '12.19'
$ws.Range("E36").Value = '  -0.98%  '
$ws.Range("D37").Value = '0.06303'
$ws.Range("E37").Value = '  +0.54%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.02317'
$ws.Range("E38").Value = '  -2.68%  '
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").Value = '5.202'
$ws.Range("E39").Value = '  -0.32%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = '0.6609'
$ws.Range("E40").Value = '  -3.19%  '
$ws.Range("D41").Value = '1.203'
$ws.Range("E41").Value = '  -1.54%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '8.064'
$ws.Range("E42").Value = '  -2.17%  '
$ws.Range("B43").Value = 'WEMIXTOKEN'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D43").Value = '1.430'
$ws.Range("E43").Value = '  -4.67%  '
$ws.Range("D44").Value = '0.9994'
$ws.Range("E44").Value = '  +0.44%  '
$ws.Range("D45").Value = '13.93'
$ws.Range("E45").Value = '  -1.24%  '
$ws.Range("D46").Value = '0.6082'
$ws.Range("E46").Value = '  -1.40%  '
$ws.Range("D47").Value = '3.796'
$ws.Range("E47").Value = '  -1.52%  '
$ws.Range("D48").Value = '127.64'
$ws.Range("E48").Value = '  -1.78%  '
$ws.Range("D49").Value = '2.023'
$ws.Range("E49").Value = '  -1.30%  '
$ws.Range("D50").Value = '0.07060'
$ws.Range("E50").Value = '  -0.81%  '
$ws.Range("D51").Value = '1.144'
$ws.Range("E51").Value = '  -2.51%  '
